$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.636.35"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.651.00"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.19"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.627"
$ws.Range("E8").Value = "  +3.89%  "
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.79"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.129.95"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.497.04"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.660.57"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.61"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.45"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.64"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.10"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.69"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.66"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "529.36"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.35"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.87"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "160.81"
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.07"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0606"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.28"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.64"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.634"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0254"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0993"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0253"
$ws.Range("E50").Value = "  +11.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.73"
$ws.Range("E51").Value = "  -1.10%  "
